$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Confirmed")
$ws.Activate()

# --- Correction: Francis A. Townsend III -> Jr. (row 29, column C) ---
$ws.Range("C29").Value = "Jr."

# --- New rows 113-119 (restored/added names) ---

# Row 113: John Carr (Dublin, Ireland - Oil & Gas)
$ws.Range("A113").Value = "John "
$ws.Range("B113").Value = "Carr"
$ws.Range("E113").Value = "Dublin"
$ws.Range("F113").Value = "Ireland"
$ws.Range("G113").Value = "Ireland"
$ws.Range("H113").Value = "Oil & Gas"

# Row 114: James Hance Jr. (N.C. - Private Equity, Carlyle Group, Operating Executive)
$ws.Range("A114").Value = "James"
$ws.Range("B114").Value = "Hance"
$ws.Range("C114").Value = "Jr."
$ws.Range("F114").Value = "N.C."
$ws.Range("G114").Value = "North Carolina"
$ws.Range("H114").Value = "Private Equity"
$ws.Range("I114").Value = "Carlyle Group"
$ws.Range("J114").Value = "Operating Executive"

# Row 115: O. Wayne Isom (New York, N.Y. - Medical, Weill Cornell Medical Center, Surgeon)
$ws.Range("A115").Value = "O. Wayne"
$ws.Range("B115").Value = "Isom"
$ws.Range("E115").Value = "New York"
$ws.Range("F115").Value = "N.Y."
$ws.Range("G115").Value = "New York"
$ws.Range("H115").Value = "Medical"
$ws.Range("I115").Value = "Weill Cornell Medical Center"
$ws.Range("J115").Value = "Surgeon"

# Row 116: David Goode (Transportation, Norfolk Southern Corp., CEO*)
$ws.Range("A116").Value = "David"
$ws.Range("B116").Value = "Goode"
$ws.Range("H116").Value = "Transportation"
$ws.Range("I116").Value = "Norfolk Southern Corp."
$ws.Range("J116").Value = "CEO*"

# Row 117: William McKnight (Augusta, Ga. - Construction, McKnight Construction Company, CEO)
$ws.Range("A117").Value = "William "
$ws.Range("B117").Value = "McKnight"
$ws.Range("E117").Value = "Augusta"
$ws.Range("F117").Value = "Ga."
$ws.Range("G117").Value = "Georgia"
$ws.Range("H117").Value = "Construction"
$ws.Range("I117").Value = "McKnight Construction Company"
$ws.Range("J117").Value = "CEO"

# Row 118: Paul Savardi (Tex. - Professional Services, Insperity, CEO)
$ws.Range("A118").Value = "Paul"
$ws.Range("B118").Value = "Savardi "
$ws.Range("F118").Value = "Tex."
$ws.Range("G118").Value = "Texas"
$ws.Range("H118").Value = "Professional Services"
$ws.Range("I118").Value = "Insperity"
$ws.Range("J118").Value = "CEO"

# Row 119: John "Jacko" Maree (Johannesburg, South Africa - Financial Services, Standard Bank Group, CEO*)
$ws.Range("A119").Value = "John ""Jacko"""
$ws.Range("B119").Value = "Maree"
$ws.Range("E119").Value = "Johannesburg"
$ws.Range("F119").Value = "South Africa"
$ws.Range("G119").Value = "South Africa"
$ws.Range("H119").Value = "Financial Services"
$ws.Range("I119").Value = "Standard Bank Group"
$ws.Range("J119").Value = "CEO*"

# --- View state: scroll to the newly-added rows and select B119 ---
[void]$ws.Range("A94").Select()
$excel.ActiveWindow.ScrollRow = 94
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B119").Select()
